# Applies the cell value updates from the 2025-02-17 FlashScore odds refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("N3").Value = 26
# Row 6
$ws.Range("Q6").Value = 3.4
$ws.Range("R6").Value = 1.33
$ws.Range("S6").Value = 7
$ws.Range("T6").Value = 1.1
$ws.Range("U6").Value = 1.75
$ws.Range("V6").Value = 2.05
# Row 10
$ws.Range("G10").Value = 1.9
$ws.Range("H10").Value = 2.88
$ws.Range("J10").Value = 2.75
$ws.Range("K10").Value = 1.8
$ws.Range("M10").Value = 1.18
$ws.Range("N10").Value = 4.5
$ws.Range("O10").Value = 1.8
$ws.Range("P10").Value = 1.91
$ws.Range("Q10").Value = 3.6
$ws.Range("R10").Value = 1.29
$ws.Range("S10").Value = 9
$ws.Range("T10").Value = 1.07
$ws.Range("U10").Value = 1.8
$ws.Range("W10").Value = 3
$ws.Range("X10").Value = 1.36
$ws.Range("Y10").Value = 4.33
$ws.Range("AE10").Value = 4.5
$ws.Range("AG10").Value = 29
$ws.Range("AH10").Value = 151
$ws.Range("AO10").Value = 81
# Row 11
$ws.Range("O11").Value = 1.67
$ws.Range("P11").Value = 2.1
$ws.Range("Y11").Value = 5
$ws.Range("Z11").Value = 8
$ws.Range("AA11").Value = 11
$ws.Range("AG11").Value = 23
# Row 13
$ws.Range("G13").Value = 4.4
$ws.Range("H13").Value = 3.65
$ws.Range("I13").Value = 1.7
$ws.Range("J13").Value = 4.6
$ws.Range("L13").Value = 2.25
$ws.Range("O13").Value = 1.24
$ws.Range("P13").Value = 3.3
$ws.Range("Q13").Value = 1.7
$ws.Range("R13").Value = 1.91
$ws.Range("S13").Value = 2.65
$ws.Range("T13").Value = 1.37
$ws.Range("W13").Value = 1.7
$ws.Range("X13").Value = 1.93
$ws.Range("Y13").Value = 13.5
$ws.Range("Z13").Value = 26
$ws.Range("AA13").Value = 14.5
$ws.Range("AB13").Value = 75
$ws.Range("AC13").Value = 40
$ws.Range("AD13").Value = 40
$ws.Range("AE13").Value = 11.5
$ws.Range("AF13").Value = 7.2
$ws.Range("AG13").Value = 14.5
$ws.Range("AH13").Value = 60
$ws.Range("AI13").Value = 450
$ws.Range("AJ13").Value = 7.6
$ws.Range("AK13").Value = 8.5
$ws.Range("AL13").Value = 8
$ws.Range("AM13").Value = 13.5
$ws.Range("AN13").Value = 13
$ws.Range("AO13").Value = 23
# Row 14
$ws.Range("G14").Value = 2
$ws.Range("I14").Value = 3.6
$ws.Range("J14").Value = 2.63
$ws.Range("Q14").Value = 1.88
$ws.Range("R14").Value = 1.98
$ws.Range("S14").Value = 3
$ws.Range("T14").Value = 1.36
$ws.Range("U14").Value = 1.36
$ws.Range("V14").Value = 3
$ws.Range("AC14").Value = 15
$ws.Range("AJ14").Value = 12
$ws.Range("AL14").Value = 13
# Row 17
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 6.5
$ws.Range("AF17").Value = 6
# Row 19
$ws.Range("G19").Value = 3.25
$ws.Range("M19").Value = 1.05
$ws.Range("N19").Value = 11
$ws.Range("O19").Value = 1.25
$ws.Range("P19").Value = 3.75
$ws.Range("Q19").Value = 1.88
$ws.Range("R19").Value = 1.98
$ws.Range("S19").Value = 3
$ws.Range("T19").Value = 1.36
$ws.Range("W19").Value = 1.67
$ws.Range("X19").Value = 2.1
$ws.Range("Y19").Value = 11
$ws.Range("AK19").Value = 11
# Row 22
$ws.Range("G22").Value = 2.18
$ws.Range("H22").Value = 2.62
$ws.Range("I22").Value = 4.2
$ws.Range("J22").Value = 2.95
$ws.Range("K22").Value = 1.75
$ws.Range("L22").Value = 5
$ws.Range("M22").Value = 1.17
$ws.Range("N22").Value = 4.4
$ws.Range("O22").Value = 1.7
$ws.Range("P22").Value = 2.05
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 1.34
$ws.Range("S22").Value = 5.6
$ws.Range("U22").Value = 1.7
$ws.Range("V22").Value = 2.05
$ws.Range("W22").Value = 2.37
$ws.Range("X22").Value = 1.52
$ws.Range("Y22").Value = 4.85
$ws.Range("Z22").Value = 8.75
$ws.Range("AA22").Value = 10
$ws.Range("AB22").Value = 22
$ws.Range("AC22").Value = 25
$ws.Range("AE22").Value = 4.4
$ws.Range("AF22").Value = 5.6
$ws.Range("AG22").Value = 22
$ws.Range("AJ22").Value = 7.5
$ws.Range("AK22").Value = 21
$ws.Range("AL22").Value = 16
$ws.Range("AM22").Value = 80
$ws.Range("AN22").Value = 65
$ws.Range("AO22").Value = 90
# Row 23
$ws.Range("G23").Value = 5.8
$ws.Range("H23").Value = 3.4
$ws.Range("I23").Value = 1.62
$ws.Range("J23").Value = 5.9
$ws.Range("K23").Value = 2.05
$ws.Range("L23").Value = 2.22
$ws.Range("M23").Value = 1.1
$ws.Range("N23").Value = 5.9
$ws.Range("O23").Value = 1.45
$ws.Range("P23").Value = 2.55
$ws.Range("Q23").Value = 2.35
$ws.Range("R23").Value = 1.53
$ws.Range("S23").Value = 4.2
$ws.Range("T23").Value = 1.19
$ws.Range("V23").Value = 2.42
$ws.Range("W23").Value = 2.3
$ws.Range("X23").Value = 1.55
$ws.Range("Y23").Value = 11.25
$ws.Range("Z23").Value = 32
$ws.Range("AA23").Value = 19.5
$ws.Range("AB23").Value = 120
$ws.Range("AC23").Value = 75
$ws.Range("AE23").Value = 5.9
$ws.Range("AF23").Value = 6.9
$ws.Range("AG23").Value = 23
$ws.Range("AJ23").Value = 4.9
$ws.Range("AK23").Value = 6.2
$ws.Range("AM23").Value = 11.5
$ws.Range("AN23").Value = 16.5
# Row 25
$ws.Range("Q25").Value = 2.05
$ws.Range("R25").Value = 1.75
# Row 28
$ws.Range("N28").Value = 15
$ws.Range("U28").Value = 1.3
$ws.Range("V28").Value = 3.4
$ws.Range("AB28").Value = 9
$ws.Range("AD28").Value = 26
$ws.Range("AO28").Value = 51
$ws.Range("AP28").Value = 2.03
$ws.Range("AQ28").Value = 1.83
# Row 29
$ws.Range("K29").Value = 2.75
$ws.Range("N29").Value = 26
$ws.Range("U29").Value = 1.2
$ws.Range("V29").Value = 4.33
$ws.Range("AI29").Value = 81
$ws.Range("AL29").Value = 15
# Row 30
$ws.Range("Q30").Value = 1.7
$ws.Range("R30").Value = 2.1
$ws.Range("S30").Value = 2.63
$ws.Range("T30").Value = 1.44
$ws.Range("AK30").Value = 12
$ws.Range("AN30").Value = 15
# Row 31
$ws.Range("G31").Value = 2.88
$ws.Range("I31").Value = 2.63
$ws.Range("J31").Value = 3.5
$ws.Range("L31").Value = 3.25
$ws.Range("AJ31").Value = 8
$ws.Range("AS31").Value = 2.28
# Row 32
$ws.Range("O32").Value = 1.22
$ws.Range("P32").Value = 4
$ws.Range("Q32").Value = 1.8
$ws.Range("R32").Value = 2
$ws.Range("S32").Value = 2.75
$ws.Range("T32").Value = 1.4
# Row 33
$ws.Range("G33").Value = 4.33
$ws.Range("H33").Value = 3.8
$ws.Range("J33").Value = 5
$ws.Range("K33").Value = 2.2
$ws.Range("L33").Value = 2.38
$ws.Range("M33").Value = 1.05
$ws.Range("N33").Value = 11
$ws.Range("O33").Value = 1.29
$ws.Range("P33").Value = 3.5
$ws.Range("Q33").Value = 1.95
$ws.Range("R33").Value = 1.9
$ws.Range("S33").Value = 3.4
$ws.Range("T33").Value = 1.3
$ws.Range("U33").Value = 1.4
$ws.Range("V33").Value = 2.75
$ws.Range("W33").Value = 1.83
$ws.Range("X33").Value = 1.83
$ws.Range("Z33").Value = 21
$ws.Range("AA33").Value = 15
$ws.Range("AE33").Value = 11
$ws.Range("AG33").Value = 17
$ws.Range("AI33").Value = 301
$ws.Range("AJ33").Value = 7
$ws.Range("AN33").Value = 15
$ws.Range("AO33").Value = 26
# Row 34
$ws.Range("Q34").Value = 2.6
$ws.Range("R34").Value = 1.48
$ws.Range("AR34").Value = 1.95
$ws.Range("AS34").Value = 1.9
# Row 39
$ws.Range("M39").Value = 1.05
$ws.Range("N39").Value = 11
$ws.Range("Q39").Value = 2.03
$ws.Range("R39").Value = 1.83
# Row 41
$ws.Range("G41").Value = 4.5
$ws.Range("H41").Value = 4.33
$ws.Range("I41").Value = 1.62
$ws.Range("J41").Value = 4.75
$ws.Range("K41").Value = 2.5
$ws.Range("L41").Value = 2.1
$ws.Range("W41").Value = 1.53
$ws.Range("X41").Value = 2.38
$ws.Range("AD41").Value = 29
$ws.Range("AE41").Value = 19
$ws.Range("AF41").Value = 9
$ws.Range("AG41").Value = 13
$ws.Range("AJ41").Value = 11
$ws.Range("AK41").Value = 10
$ws.Range("AN41").Value = 12
# Row 42
$ws.Range("J42").Value = 2.88
$ws.Range("K42").Value = 2.25
$ws.Range("AK42").Value = 17
# Row 43
$ws.Range("L43").Value = 3.25
$ws.Range("M43").Value = 1.06
$ws.Range("N43").Value = 10
$ws.Range("U43").Value = 1.44
$ws.Range("V43").Value = 2.63
$ws.Range("AM43").Value = 26
# Row 44
$ws.Range("M44").Value = 1.05
$ws.Range("N44").Value = 11
# Row 45
$ws.Range("G45").Value = 2.45
$ws.Range("I45").Value = 2.8
$ws.Range("O45").Value = 1.29
$ws.Range("P45").Value = 3.5
$ws.Range("Q45").Value = 1.95
$ws.Range("R45").Value = 1.85
$ws.Range("S45").Value = 3.4
$ws.Range("T45").Value = 1.3
$ws.Range("U45").Value = 1.4
$ws.Range("V45").Value = 2.75
$ws.Range("W45").Value = 1.73
$ws.Range("X45").Value = 2
$ws.Range("Y45").Value = 9
$ws.Range("AE45").Value = 10
$ws.Range("AG45").Value = 13
$ws.Range("AH45").Value = 41
$ws.Range("AI45").Value = 201
$ws.Range("AJ45").Value = 9.5
$ws.Range("AN45").Value = 21
$ws.Range("AO45").Value = 29
